# Updated symbol list on Wed Jan 25 20:56:57 UTC 2023 with GitHub Actions
#
# The "Price" (column D) and "Volume(1h)" (column E) figures on the coin
# tracker sheet are refreshed with the latest scrape. Every value in these
# two columns is stored as literal text (e.g. "303.03", "-2.31%") rather
# than as a number/percentage, so writing it back with a plain
# `Range.Value =` assignment would let Excel auto-coerce numeric-looking
# strings into real numbers (and "%"-suffixed strings into percentage
# numbers), which would silently reformat the cell and corrupt the
# intended plain-text representation. To avoid that, each cell is
# temporarily marked as Text ("@") before the new literal is written, and
# the style is reset back to Normal immediately after so no stray
# number-format/style ends up attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "303.03"
    "E2" = "-2.31%"
    "D3" = "35.97"
    "E3" = "1.22%"
    "D4" = "5.071"
    "E4" = "-0.95%"
    "D5" = "0.08066"
    "E5" = "-1.58%"
    "D6" = "1.962"
    "E6" = "-4.07%"
    "D7" = "7.800"
    "E7" = "-2.06%"
    "D8" = "0.9284"
    "E8" = "0.16%"
    "D9" = "0.1499"
    "E9" = "38.01%"
    "D10" = "0.1896"
    "E10" = "-1.42%"
    "D11" = "0.08999"
    "E11" = "-4.60%"
    "D12" = "0.03455"
    "E12" = "-4.22%"
    "D13" = "0.09847"
    "E13" = "-0.39%"
    "D14" = "0.001392"
    "E14" = "-2.73%"
    "D15" = "0.005780"
    "E15" = "1.05%"
    "D16" = "3.538"
    "E16" = "2.04%"
    "D17" = "4.050"
    "E17" = "-1.95%"
    "D18" = "2.962"
    "E18" = "-0.96%"
    "D19" = "0.3443"
    "E19" = "0.72%"
    "D20" = "0.1298"
    "D21" = "5.033"
    "E21" = "-1.11%"
    "D22" = "0.2391"
    "E22" = "9.06%"
    "D23" = "0.04498"
    "E23" = "-1.04%"
    "D24" = "0.001204"
    "E24" = "-1.66%"
    "D25" = "0.004811"
    "E25" = "0.57%"
    "D26" = "0.0001225"
    "E26" = "-1.96%"
    "E27" = "-32.28%"
    "D39" = "0.01880"
    "E39" = "-4.70%"
    "E40" = "-2.06%"
    "D41" = "0.01059"
    "E41" = "7.66%"
    "D42" = "0.007335"
    "E42" = "-5.90%"
    "E43" = "-2.65%"
    "D44" = "0.002102"
    "E44" = "-0.61%"
    "D45" = "0.009727"
    "E45" = "-15.80%"
    "D46" = "0.00006209"
    "E46" = "-4.54%"
    "D47" = "0.00000000747"
    "E47" = "-0.44%"
    "E48" = "-62.62%"
    "D50" = "0.00002093"
    "E50" = "-0.44%"
    "D51" = "0.0001993"
    "E51" = "-0.44%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
